$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a value to be stored as text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "26.100.05", "1.000", "0.06198") into numbers.
function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Coin name / link swap (rows 50 and 51) ---
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# --- Price column (column D) updates ---
Set-TextValue 'D2' '26.100.05'
Set-TextValue 'D3' '1.751.22'
Set-TextValue 'D4' '1.000'
Set-TextValue 'D5' '235.39'
Set-TextValue 'D6' '0.9997'
Set-TextValue 'D7' '0.5297'
Set-TextValue 'D9' '0.06198'
Set-TextValue 'D10' '1.744.33'
Set-TextValue 'D11' '0.07178'
Set-TextValue 'D12' '15.48'
Set-TextValue 'D13' '0.6480'
Set-TextValue 'D14' '4.634'
Set-TextValue 'D15' '78.78'
Set-TextValue 'D17' '1.000'
Set-TextValue 'D18' '26.001.91'
Set-TextValue 'D20' '0.000006747'
Set-TextValue 'D21' '1.968.48'
Set-TextValue 'D22' '4.328'
Set-TextValue 'D23' '8.746'
Set-TextValue 'D24' '5.256'
Set-TextValue 'D25' '138.97'
Set-TextValue 'D26' '1.514'
Set-TextValue 'D27' '15.32'
Set-TextValue 'D28' '1.814'
Set-TextValue 'D29' '104.63'
Set-TextValue 'D30' '0.08314'
Set-TextValue 'D31' '3.824'
Set-TextValue 'D32' '3.663'
Set-TextValue 'D33' '0.04591'
Set-TextValue 'D35' '1.012'
Set-TextValue 'D36' '0.6368'
Set-TextValue 'D37' '2.715'
Set-TextValue 'D38' '0.01607'
Set-TextValue 'D39' '1.965'
Set-TextValue 'D40' '0.9991'
Set-TextValue 'D41' '100.72'
Set-TextValue 'D42' '0.3950'
Set-TextValue 'D43' '0.7481'
Set-TextValue 'D44' '5.040'
Set-TextValue 'D45' '0.1152'
Set-TextValue 'D46' '6.383'
Set-TextValue 'D47' '0.05351'
Set-TextValue 'D48' '31.13'
Set-TextValue 'D49' '54.45'
Set-TextValue 'D50' '0.3472'
Set-TextValue 'D51' '7.610'

# --- Volume(1h) column (column E) updates ---
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E5').Value = '  +4.41%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  +2.59%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('E22').Value = '  +6.35%  '
$ws.Range('E23').Value = '  +3.84%  '
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  +2.47%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('E31').Value = '  +6.09%  '
$ws.Range('E32').Value = '  +7.63%  '
$ws.Range('E33').Value = '  +4.54%  '
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('E36').Value = '  +6.20%  '
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('E38').Value = '  +3.88%  '
$ws.Range('E39').Value = '  +3.53%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  +3.47%  '
$ws.Range('E43').Value = '  +3.42%  '
$ws.Range('E44').Value = '  +3.24%  '
$ws.Range('E45').Value = '  +4.95%  '
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('E48').Value = '  +4.82%  '
$ws.Range('E49').Value = '  +4.52%  '
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('E51').Value = '  +1.65%  '
